# ------------------------------------------------------------------
# Adds a new "2022-Q3" sheet (with fund-holding detail) right after
# the "总计" summary sheet, and updates the summary sheet with a new
# row of aggregated 2022-Q3 data.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# 1) Create the new "2022-Q3" worksheet, inserted right before the
#    sheet that is currently in position 2 (the old "2022-Q2" sheet).
# ====================================================================
$summarySheet = $wb.Worksheets.Item(1)
$oldSecondSheet = $wb.Worksheets.Item(2)

$q3 = $wb.Worksheets.Add($oldSecondSheet)
$q3.Name = "2022-Q3"

# --- header row (row 1), columns B..H ------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $cell = $q3.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- fund holdings data (rows 2..11) --------------------------------
# columns: code, name, size, stock-position, ratio, held-value, rank
$rows = @(
    @("161810", "银华内需精选混合（LOF）", "23.47", "94.62", "5.87", "1.3777", 9),
    @("004475", "华泰柏瑞富利灵活配置混合A", "37.71", "67.54", "1.65", "0.6222", 10),
    @("014597", "华泰柏瑞富利灵活配置混合C", "21.79", "67.54", "1.65", "0.3595", 10),
    @("180020", "银华成长先锋混合", "2.15", "79.28", "5.26", "0.1131", 8),
    @("011429", "前海开源民裕进取混合", "2.53", "60.79", "3.10", "0.0784", 6),
    @("008480", "永赢股息优选混合A", "2.05", "67.89", "3.06", "0.0627", 9),
    @("003175", "华泰柏瑞多策略灵活配置混合A", "3.32", "67.70", "1.76", "0.0584", 9),
    @("011588", "前海开源成份精选混合", "0.84", "60.90", "3.62", "0.0304", 4),
    @("015450", "华泰柏瑞多策略灵活配置混合C", "1.15", "67.70", "1.76", "0.0202", 9),
    @("008481", "永赢股息优选混合C", "0.17", "67.89", "3.06", "0.0052", 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $rowData = $rows[$r]

    # column A: numeric index, bold/centered/bordered like the header
    $aCell = $q3.Cells.Item($rowNum, 1)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # columns B..G: text values (force text so numeric-looking strings
    # such as "23.47" are NOT converted into real numbers)
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q3.Cells.Item($rowNum, 2 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
    }

    # column H: numeric rank
    $q3.Cells.Item($rowNum, 8).Value = $rowData[6]
}

# ====================================================================
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q3
#    right after the header row, pushing the rest down, then keep the
#    running index in column A sequential (0..7).
# ====================================================================
$summarySheet.Rows(2).Insert()

# copy formatting so the new row looks like the existing rows
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summarySheet.Range("B3:D3").Copy()
$summarySheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("C2").Value = 10
$summarySheet.Range("D2").Value = 2.73

$summarySheet.Range("A3").Value = 1
$summarySheet.Range("A4").Value = 2
$summarySheet.Range("A5").Value = 3
$summarySheet.Range("A6").Value = 4
$summarySheet.Range("A7").Value = 5
$summarySheet.Range("A8").Value = 6
$summarySheet.Range("A9").Value = 7

Write-Host "2022-Q3 sheet added and summary sheet updated."
